$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New benchmark rows appended below existing data (row 6 and row 7)
$ws.Range("A7").Value = "Func"
$ws.Range("B7").Value = "Pa"
$ws.Range("A6").Value = "Including type conversion"
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = "Object { size: 250410, visits: 2455, time: 2565.622217476026 }"
$ws.Range("E7").Value = "Object { size: 46512, visits: 2448, time: 3946.8605084212986 }"

# Update view/selection state to match target (also clears the stale topLeftCell scroll anchor)
$ws.Range("E7").Select()
